$wb = $excel.ActiveWorkbook

# Remove the obsolete "Sheet" row (row 16) from the optimization_parameters sheet.
# This row contained the label "Sheet" in column A and the values 3/4 in B/C,
# which are no longer used; deleting it shifts the simulation_timepoints row
# up from row 17 to row 16.
$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows.Item(16).Delete()

# Make "optimization_diagnostics" the active sheet/tab.
$diagSheet = $wb.Worksheets.Item("optimization_diagnostics")
$diagSheet.Activate()
